# San Antonio roster: Blake Wesley and Gorgui Dieng swap roster rows (row 12 <-> row 13).
# Column A ("No." index column, 10/11) stays put; all other columns (B:K) - jersey #, player,
# position, height, weight, birth date, country code, experience, college, bbref url - move
# together with the player, so we swap the full B:K content between row 12 and row 13.
#
# We use Copy + PasteSpecial(xlPasteValues) through a scratch cell rather than plain
# Range.Value assignment so that text-typed numeric-looking cells (e.g. the "Exp" column
# storing "9" as text) keep their original text type instead of being coerced to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

$srcRow = 12
$dstRow = 13
$colFirst = "B"
$colLast = "K"
$scratch = "Z1"

$rangeSrc = "$colFirst$srcRow" + ":" + "$colLast$srcRow"
$rangeDst = "$colFirst$dstRow" + ":" + "$colLast$dstRow"

# 1) Stash row 12's current values in a scratch area.
$ws.Range($rangeSrc).Copy() | Out-Null
$ws.Range($scratch).PasteSpecial($xlPasteValues) | Out-Null

# 2) Move row 13's values into row 12.
$ws.Range($rangeDst).Copy() | Out-Null
$ws.Range($rangeSrc).PasteSpecial($xlPasteValues) | Out-Null

# 3) Move the stashed original row 12 values into row 13.
$scratchRange = $ws.Range($scratch).Resize(1, 10)
$scratchRange.Copy() | Out-Null
$ws.Range($rangeDst).PasteSpecial($xlPasteValues) | Out-Null

# 4) Clean up the scratch area.
$scratchRange.Clear() | Out-Null

$wb.Application.CutCopyMode = $false
